# The 'ObjTables' SBtab schema metadata is stored as literal header text in row 1
# (and, for the first sheet, row 1 holds a workbook-wide banner with row 2 holding
# the per-table header) of every worksheet. Each header encodes the objTablesVersion
# and a generation timestamp. Bump objTablesVersion 0.0.9 -> 1.0.0 and refresh the
# date 2020-04-27 01:05:49 -> 2020-05-29 00:19:44 across every sheet, preserving the
# rest of each header line (class/name/document attributes) verbatim.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('!!Compartment')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!!ObjTables schema=''SBtab'' objTablesVersion=''1.0.0'' date=''2020-05-29 00:19:44'''
$ws.Range("A2").Locked = $false
$ws.Range("A2").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Compartment'' name=''Compartment'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Compound')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Compound'' name=''Compound'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Definition')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Definition'' name=''Definition'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Enzyme')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Enzyme'' name=''Enzyme'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!FbcObjective')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''FbcObjective'' name=''FbcObjective'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Gene')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Gene'' name=''Gene'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'' document=''lac_Operon'''

$ws = $wb.Worksheets.Item('!!Layout')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Layout'' name=''Layout'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Measurement')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Measurement'' name=''Measurement'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!PbConfig')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''PbConfig'' name=''PbConfig'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Position')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Position'' name=''Position'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Protein')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Protein'' name=''Protein'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Quantity')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Quantity'' name=''Quantity'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!QuantityInfo')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''QuantityInfo'' name=''QuantityInfo'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!QuantityMatrix')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Reaction')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Reaction'' name=''Reaction'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!ReactionStoichiometry')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Regulator')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Regulator'' name=''Regulator'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'' document=''lac_Operon'''

$ws = $wb.Worksheets.Item('!!Relation')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Relation'' name=''Relation'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!Relationship')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''Relationship'' name=''Relationship'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!SparseMatrix')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''SparseMatrix'' name=''SparseMatrix'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!SparseMatrixColumn')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!SparseMatrixOrdered')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!SparseMatrixRow')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!StoichiometricMatrix')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!rxnconContingencyList')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''

$ws = $wb.Worksheets.Item('!!rxnconReactionList')
$ws.Range("A1").Locked = $false
$ws.Range("A1").Value = '!!ObjTables schema=''SBtab'' type=''Data'' tableFormat=''row'' class=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-05-29 00:19:44'' objTablesVersion=''1.0.0'''
